# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# J.Taylor (row 3)
$rushing.Cells.Item(3, 3).Value = 1
$rushing.Cells.Item(3, 4).Value = 3

# N.Hines (row 4)
$rushing.Cells.Item(4, 3).Value = 175
$rushing.Cells.Item(4, 4).Value = 119
$rushing.Cells.Item(4, 5).Value = 37
$rushing.Cells.Item(4, 6).Value = 82

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# J.Taylor (row 2)
$receiving.Cells.Item(2, 3).Value = 50
$receiving.Cells.Item(2, 4).Value = 37
$receiving.Cells.Item(2, 7).Value = 6
$receiving.Cells.Item(2, 8).Value = 5

# N.Hines (row 3)
$receiving.Cells.Item(3, 3).Value = 50
$receiving.Cells.Item(3, 4).Value = 36
$receiving.Cells.Item(3, 7).Value = 3
$receiving.Cells.Item(3, 8).Value = 3

# M.Pittman (row 5)
$receiving.Cells.Item(5, 3).Value = 104
$receiving.Cells.Item(5, 4).Value = 77
$receiving.Cells.Item(5, 5).Value = 25
$receiving.Cells.Item(5, 6).Value = 12
$receiving.Cells.Item(5, 7).Value = 18
$receiving.Cells.Item(5, 8).Value = 9

# Z.Pascal (row 6)
$receiving.Cells.Item(6, 3).Value = 55
$receiving.Cells.Item(6, 4).Value = 31
$receiving.Cells.Item(6, 7).Value = 13
$receiving.Cells.Item(6, 8).Value = 6

# A.Dulin (row 7)
$receiving.Cells.Item(7, 3).Value = 13
$receiving.Cells.Item(7, 5).Value = 7

# M.Strachan (row 10)
$receiving.Cells.Item(10, 3).Value = 41
$receiving.Cells.Item(10, 4).Value = 31
$receiving.Cells.Item(10, 5).Value = 10
$receiving.Cells.Item(10, 6).Value = 6
$receiving.Cells.Item(10, 7).Value = 4

# M.Alie-Cox (row 13)
$receiving.Cells.Item(13, 3).Value = 39
$receiving.Cells.Item(13, 4).Value = 27

# K.Granson (row 14)
$receiving.Cells.Item(14, 3).Value = 34
$receiving.Cells.Item(14, 4).Value = 21
$receiving.Cells.Item(14, 5).Value = 14
